# Add new "PRT Overwrite" option row (row 33) to the QualityChecks template.
# Mirrors the existing "MTN.OVERWRITE" style row: a Value label, a FALSE
# checkbox, a description, and a Field_ID code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row 33 cells. Shared strings get appended in B, E, D
# order so the resulting sharedStrings.xml index assignment matches the
# target file (84=label, 85=field id, 86=description).
$ws.Range("B33").Value = "PRT Overwrite"
$ws.Range("C33").Value = $false
$ws.Range("E33").Value = "PRT.OVERWRITE"
$ws.Range("D33").Value = "When copying PRT, overwrite existing files."

# Move the view/selection down to the new row, matching the author's
# on-save cursor position.
$ws.Range("E33").Select() | Out-Null
